# Move the "_GoBack" bookmark from the end of the "امیرمزلقانی" paragraph
# down to the (now) final empty paragraph of the document, and center
# that final paragraph.

$d = $word.ActiveDocument

# --- 1. Remove the existing "_GoBack" bookmark -----------------------
# "_GoBack" is a hidden bookmark (its name starts with "_"), so it will
# never show up in $d.Bookmarks / enumeration or Count, but it can still
# be reached directly by name.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Center the last paragraph of the document body ---------------
# wdAlignParagraphCenter = 1
$lastParaIndex = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($lastParaIndex)
$targetPara.Format.Alignment = 1

# --- 3. Re-create "_GoBack" on that last (empty) paragraph ------------
# Building the bookmark range from the paragraph's own Range object
# directly does not anchor correctly when the paragraph has no runs
# (it is empty), so the range is rebuilt from plain numeric offsets,
# starting one character before the paragraph so the insertion point
# resolves against real content.
$paraStart = $targetPara.Range.Start
$paraEnd = $targetPara.Range.End
$bookmarkRange = $d.Range($paraStart - 1, $paraEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
